$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B7 / F7 / B12 from "time" to "special"
$ws.Range("B7").Value = "special"
$ws.Range("F7").Value = "special"
$ws.Range("B12").Value = "special"

# Update the selection to B7 (single cell, matches the diff's <selection activeCell="B7" sqref="B7"/>)
$ws.Range("B7").Select()
